$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 2 (header row) - the [h,k,l] combo columns were re-ordered
#    while the 1Pair/2Pairs/.../MaxUnique columns stayed the same.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 3).Value  = "[5, 1, 1]"
$ws.Cells.Item(2, 4).Value  = "[4, 2, 2]"
$ws.Cells.Item(2, 5).Value  = "[3, 1, 1]"
$ws.Cells.Item(2, 6).Value  = "[3, 3, 1]"
$ws.Cells.Item(2, 7).Value  = "[2, 2, 2]"
$ws.Cells.Item(2, 8).Value  = "[1, 1, 1]"
$ws.Cells.Item(2, 9).Value  = "[3, 3, 3]"
$ws.Cells.Item(2, 10).Value = "[2, 2, 0]"
$ws.Cells.Item(2, 11).Value = "[2, 0, 0]"
$ws.Cells.Item(2, 12).Value = "[4, 0, 0]"
$ws.Cells.Item(2, 13).Value = "[4, 2, 0]"

# ---------------------------------------------------------------------------
# 2. Rename the existing simulation rows (rows 3-19, column B) to the new
#    scheme names used in the updated report.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 2).Value  = "Spiral5"
$ws.Cells.Item(4, 2).Value  = "RotRing OmegaMax-90"
$ws.Cells.Item(5, 2).Value  = "Equal Angle"
$ws.Cells.Item(6, 2).Value  = "Tilt Rotate"
$ws.Cells.Item(7, 2).Value  = "CLR"
$ws.Cells.Item(8, 2).Value  = "Rizzie Hex"
$ws.Cells.Item(9, 2).Value  = "Thomas Hex"
$ws.Cells.Item(10, 2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(11, 2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(12, 2).Value = "Equal Angle_Partial"
$ws.Cells.Item(13, 2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(14, 2).Value = "ND Single"
$ws.Cells.Item(15, 2).Value = "RD Single"
$ws.Cells.Item(16, 2).Value = "TD Single"
$ws.Cells.Item(17, 2).Value = "Morris Single"
$ws.Cells.Item(18, 2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(19, 2).Value = "Ring Perpendicular to RD"

# ---------------------------------------------------------------------------
# 3. Add new simulation rows 20-29, copying formatting from row 19 first so
#    that column A keeps the same bold/bordered style.
# ---------------------------------------------------------------------------
$ws.Range("A19:W19").Copy()
$ws.Range("A20:A29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ A = 18; B = "Ring Perpendicular to TD" },
    @{ A = 19; B = "OffsetFTD" },
    @{ A = 20; B = "OffsetATD" },
    @{ A = 21; B = "OffsetF45" },
    @{ A = 22; B = "OffsetA45" },
    @{ A = 23; B = "OffsetFRD" },
    @{ A = 24; B = "OffsetARD" },
    @{ A = 25; B = "Gaussian Quadrature" },
    @{ A = 26; B = "Michael-CCHex" },
    @{ A = 27; B = "Michael-SNHex" }
)

$rowNum = 20
foreach ($item in $newRows) {
    $ws.Cells.Item($rowNum, 1).Value = $item.A
    $ws.Cells.Item($rowNum, 2).Value = $item.B
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($rowNum, $c).Value = 1
    }
    $rowNum = $rowNum + 1
}
